$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("three_line")
$ws.Cells.Item(367, 1).Value = 45306
$ws.Cells.Item(367, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(367, 2).Value = "17-06-2024 00:00:00"
$ws.Cells.Item(367, 3).Value = "week"
$ws.Cells.Item(367, 4).Value = "DMART.NS"
$ws.Cells.Item(367, 5).Value = 45096
$ws.Cells.Item(367, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(367, 6).Value = 4150
$ws.Cells.Item(367, 7).Value = 45264
$ws.Cells.Item(367, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(367, 8).Value = 4206
$ws.Cells.Item(367, 9).Value = 0
$ws.Cells.Item(367, 10).Value = 0
$ws.Cells.Item(367, 11).Value = "High"
$ws.Cells.Item(367, 12).Value = "18/06/2024 11:36:33"

$ws.Cells.Item(368, 1).Value = 45362
$ws.Cells.Item(368, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(368, 2).Value = "17-06-2024 00:00:00"
$ws.Cells.Item(368, 3).Value = "week"
$ws.Cells.Item(368, 4).Value = "MAHEPC.NS"
$ws.Cells.Item(368, 5).Value = 44928
$ws.Cells.Item(368, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(368, 6).Value = 113
$ws.Cells.Item(368, 7).Value = 45320
$ws.Cells.Item(368, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(368, 8).Value = 162.8000030517578
$ws.Cells.Item(368, 9).Value = 0
$ws.Cells.Item(368, 10).Value = 0
$ws.Cells.Item(368, 11).Value = "High"
$ws.Cells.Item(368, 12).Value = "18/06/2024 11:36:33"

$ws.Cells.Item(369, 1).Value = 45446
$ws.Cells.Item(369, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(369, 2).Value = "17-06-2024 00:00:00"
$ws.Cells.Item(369, 3).Value = "week"
$ws.Cells.Item(369, 4).Value = "GMRINFRA.NS"
$ws.Cells.Item(369, 5).Value = 45341
$ws.Cells.Item(369, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(369, 6).Value = 94.34999847412109
$ws.Cells.Item(369, 7).Value = 45404
$ws.Cells.Item(369, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(369, 8).Value = 92.40000152587891
$ws.Cells.Item(369, 9).Value = 0
$ws.Cells.Item(369, 10).Value = 0
$ws.Cells.Item(369, 11).Value = "High"
$ws.Cells.Item(369, 12).Value = "18/06/2024 11:36:33"

$ws.Cells.Item(370, 1).Value = 45369
$ws.Cells.Item(370, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(370, 2).Value = "17-06-2024 00:00:00"
$ws.Cells.Item(370, 3).Value = "week"
$ws.Cells.Item(370, 4).Value = "JISLJALEQS.NS"
$ws.Cells.Item(370, 5).Value = 45250
$ws.Cells.Item(370, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(370, 6).Value = 73.44999694824219
$ws.Cells.Item(370, 7).Value = 45327
$ws.Cells.Item(370, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(370, 8).Value = 68.59999847412109
$ws.Cells.Item(370, 9).Value = 0
$ws.Cells.Item(370, 10).Value = 0
$ws.Cells.Item(370, 11).Value = "High"
$ws.Cells.Item(370, 12).Value = "18/06/2024 11:36:33"

$ws.Cells.Item(371, 1).Value = 45369
$ws.Cells.Item(371, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(371, 2).Value = "17-06-2024 00:00:00"
$ws.Cells.Item(371, 3).Value = "week"
$ws.Cells.Item(371, 4).Value = "M&MFIN.NS"
$ws.Cells.Item(371, 5).Value = 45271
$ws.Cells.Item(371, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(371, 6).Value = 296
$ws.Cells.Item(371, 7).Value = 45327
$ws.Cells.Item(371, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(371, 8).Value = 302.8999938964844
$ws.Cells.Item(371, 9).Value = 0
$ws.Cells.Item(371, 10).Value = 0
$ws.Cells.Item(371, 11).Value = "High"
$ws.Cells.Item(371, 12).Value = "18/06/2024 11:36:33"

$ws.Cells.Item(372, 1).Value = 45362
$ws.Cells.Item(372, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(372, 2).Value = "17-06-2024 00:00:00"
$ws.Cells.Item(372, 3).Value = "week"
$ws.Cells.Item(372, 4).Value = "JKPAPER.NS"
$ws.Cells.Item(372, 5).Value = 44788
$ws.Cells.Item(372, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(372, 6).Value = 449.9500122070312
$ws.Cells.Item(372, 7).Value = 44914
$ws.Cells.Item(372, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(372, 8).Value = 453.2000122070312
$ws.Cells.Item(372, 9).Value = 0
$ws.Cells.Item(372, 10).Value = 0
$ws.Cells.Item(372, 11).Value = "High"
$ws.Cells.Item(372, 12).Value = "18/06/2024 11:36:33"

$ws.Cells.Item(373, 1).Value = 45362
$ws.Cells.Item(373, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(373, 2).Value = "17-06-2024 00:00:00"
$ws.Cells.Item(373, 3).Value = "week"
$ws.Cells.Item(373, 4).Value = "JKPAPER.NS"
$ws.Cells.Item(373, 5).Value = 44788
$ws.Cells.Item(373, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(373, 6).Value = 449.9500122070312
$ws.Cells.Item(373, 7).Value = 45320
$ws.Cells.Item(373, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(373, 8).Value = 452
$ws.Cells.Item(373, 9).Value = 0
$ws.Cells.Item(373, 10).Value = 0
$ws.Cells.Item(373, 11).Value = "High"
$ws.Cells.Item(373, 12).Value = "18/06/2024 11:36:33"

$ws.Cells.Item(374, 1).Value = 45362
$ws.Cells.Item(374, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(374, 2).Value = "17-06-2024 00:00:00"
$ws.Cells.Item(374, 3).Value = "week"
$ws.Cells.Item(374, 4).Value = "JKPAPER.NS"
$ws.Cells.Item(374, 5).Value = 44788
$ws.Cells.Item(374, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(374, 6).Value = 449.9500122070312
$ws.Cells.Item(374, 7).Value = 45215
$ws.Cells.Item(374, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(374, 8).Value = 415.2999877929688
$ws.Cells.Item(374, 9).Value = 0
$ws.Cells.Item(374, 10).Value = 0
$ws.Cells.Item(374, 11).Value = "High"
$ws.Cells.Item(374, 12).Value = "18/06/2024 11:36:33"

$ws.Cells.Item(375, 1).Value = 44319
$ws.Cells.Item(375, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(375, 2).Value = "17-06-2024 00:00:00"
$ws.Cells.Item(375, 3).Value = "week"
$ws.Cells.Item(375, 4).Value = "BLUECLOUDS.BO"
$ws.Cells.Item(375, 5).Value = 43549
$ws.Cells.Item(375, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(375, 6).Value = 13.10999965667725
$ws.Cells.Item(375, 7).Value = 44200
$ws.Cells.Item(375, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(375, 8).Value = 12.52000045776367
$ws.Cells.Item(375, 9).Value = 0
$ws.Cells.Item(375, 10).Value = 0
$ws.Cells.Item(375, 11).Value = "High"
$ws.Cells.Item(375, 12).Value = "18/06/2024 11:36:33"

$ws.Cells.Item(376, 1).Value = 44319
$ws.Cells.Item(376, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(376, 2).Value = "17-06-2024 00:00:00"
$ws.Cells.Item(376, 3).Value = "week"
$ws.Cells.Item(376, 4).Value = "BLUECLOUDS.BO"
$ws.Cells.Item(376, 5).Value = 44200
$ws.Cells.Item(376, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(376, 6).Value = 12.52000045776367
$ws.Cells.Item(376, 7).Value = 44277
$ws.Cells.Item(376, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(376, 8).Value = 10.97000026702881
$ws.Cells.Item(376, 9).Value = 0
$ws.Cells.Item(376, 10).Value = 0
$ws.Cells.Item(376, 11).Value = "High"
$ws.Cells.Item(376, 12).Value = "18/06/2024 11:36:33"

$ws.Cells.Item(377, 1).Value = 44319
$ws.Cells.Item(377, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(377, 2).Value = "17-06-2024 00:00:00"
$ws.Cells.Item(377, 3).Value = "week"
$ws.Cells.Item(377, 4).Value = "BLUECLOUDS.BO"
$ws.Cells.Item(377, 5).Value = 43542
$ws.Cells.Item(377, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(377, 6).Value = 13.10999965667725
$ws.Cells.Item(377, 7).Value = 44200
$ws.Cells.Item(377, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(377, 8).Value = 12.52000045776367
$ws.Cells.Item(377, 9).Value = 0
$ws.Cells.Item(377, 10).Value = 0
$ws.Cells.Item(377, 11).Value = "High"
$ws.Cells.Item(377, 12).Value = "18/06/2024 11:36:33"

$ws.Cells.Item(378, 1).Value = 44319
$ws.Cells.Item(378, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(378, 2).Value = "17-06-2024 00:00:00"
$ws.Cells.Item(378, 3).Value = "week"
$ws.Cells.Item(378, 4).Value = "BLUECLOUDS.BO"
$ws.Cells.Item(378, 5).Value = 43542
$ws.Cells.Item(378, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(378, 6).Value = 13.10999965667725
$ws.Cells.Item(378, 7).Value = 44158
$ws.Cells.Item(378, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(378, 8).Value = 12
$ws.Cells.Item(378, 9).Value = 0
$ws.Cells.Item(378, 10).Value = 0
$ws.Cells.Item(378, 11).Value = "High"
$ws.Cells.Item(378, 12).Value = "18/06/2024 11:36:33"


$ws = $wb.Worksheets.Item("two_line")
$ws.Cells.Item(207, 1).Value = 45215
$ws.Cells.Item(207, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(207, 2).Value = "17-06-2024 00:00:00"
$ws.Cells.Item(207, 3).Value = "week"
$ws.Cells.Item(207, 4).Value = "BANARISUG.NS"
$ws.Cells.Item(207, 5).Value = 45089
$ws.Cells.Item(207, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(207, 6).Value = 2930.050048828125
$ws.Cells.Item(207, 7).Value = 45173
$ws.Cells.Item(207, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(207, 8).Value = 2930.300048828125
$ws.Cells.Item(207, 9).Value = "High"
$ws.Cells.Item(207, 10).Value = "18/06/2024 11:36:33"

$ws.Cells.Item(208, 1).Value = 44585
$ws.Cells.Item(208, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(208, 2).Value = "17-06-2024 00:00:00"
$ws.Cells.Item(208, 3).Value = "week"
$ws.Cells.Item(208, 4).Value = "IDEA.NS"
$ws.Cells.Item(208, 5).Value = 44536
$ws.Cells.Item(208, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(208, 6).Value = 16.79999923706055
$ws.Cells.Item(208, 7).Value = 44543
$ws.Cells.Item(208, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(208, 8).Value = 16.79999923706055
$ws.Cells.Item(208, 9).Value = "High"
$ws.Cells.Item(208, 10).Value = "18/06/2024 11:36:33"


$ws = $wb.Worksheets.Item("ph_pl_breakout_line")
$ws.Cells.Item(1329, 1).Value = "CHAMBLFERT.NS"
$ws.Cells.Item(1329, 2).Value = 45411
$ws.Cells.Item(1329, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1329, 3).Value = 439
$ws.Cells.Item(1329, 4).Value = 400.3500061035156
$ws.Cells.Item(1329, 5).Value = 405.7000122070312
$ws.Cells.Item(1329, 6).Value = "High"
$ws.Cells.Item(1329, 7).Value = 439
$ws.Cells.Item(1329, 8).Value = "week"
$ws.Cells.Item(1329, 9).Value = "17-06-2024 00:00:00"
$ws.Cells.Item(1329, 10).Value = 439.8999938964844
$ws.Cells.Item(1329, 11).Value = 437
$ws.Cells.Item(1329, 12).Value = "18/06/2024 11:36:33"

$ws.Cells.Item(1330, 1).Value = "BAJFINANCE.NS"
$ws.Cells.Item(1330, 2).Value = 45404
$ws.Cells.Item(1330, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1330, 3).Value = 7419.4501953125
$ws.Cells.Item(1330, 4).Value = 6691.5
$ws.Cells.Item(1330, 5).Value = 6731.2001953125
$ws.Cells.Item(1330, 6).Value = "High"
$ws.Cells.Item(1330, 7).Value = 7419.4501953125
$ws.Cells.Item(1330, 8).Value = "week"
$ws.Cells.Item(1330, 9).Value = "17-06-2024 00:00:00"
$ws.Cells.Item(1330, 10).Value = 7429.4501953125
$ws.Cells.Item(1330, 11).Value = 7370
$ws.Cells.Item(1330, 12).Value = "18/06/2024 11:36:33"

$ws.Cells.Item(1331, 1).Value = "PATELENG.NS"
$ws.Cells.Item(1331, 2).Value = 45271
$ws.Cells.Item(1331, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1331, 3).Value = 69.69999694824219
$ws.Cells.Item(1331, 4).Value = 62.5
$ws.Cells.Item(1331, 5).Value = 63
$ws.Cells.Item(1331, 6).Value = "High"
$ws.Cells.Item(1331, 7).Value = 69.69999694824219
$ws.Cells.Item(1331, 8).Value = "week"
$ws.Cells.Item(1331, 9).Value = "17-06-2024 00:00:00"
$ws.Cells.Item(1331, 10).Value = 72.19999694824219
$ws.Cells.Item(1331, 11).Value = 69.19999694824219
$ws.Cells.Item(1331, 12).Value = "18/06/2024 11:36:33"

$ws.Cells.Item(1332, 1).Value = "BAJAJCON.NS"
$ws.Cells.Item(1332, 2).Value = 45418
$ws.Cells.Item(1332, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1332, 3).Value = 271
$ws.Cells.Item(1332, 4).Value = 236
$ws.Cells.Item(1332, 5).Value = 244.9499969482422
$ws.Cells.Item(1332, 6).Value = "High"
$ws.Cells.Item(1332, 7).Value = 271
$ws.Cells.Item(1332, 8).Value = "week"
$ws.Cells.Item(1332, 9).Value = "17-06-2024 00:00:00"
$ws.Cells.Item(1332, 10).Value = 273.8399963378906
$ws.Cells.Item(1332, 11).Value = 268
$ws.Cells.Item(1332, 12).Value = "18/06/2024 11:36:33"

$ws.Cells.Item(1333, 1).Value = "MAHLIFE.NS"
$ws.Cells.Item(1333, 2).Value = 45327
$ws.Cells.Item(1333, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1333, 3).Value = 632.7999877929688
$ws.Cells.Item(1333, 4).Value = 571.5
$ws.Cells.Item(1333, 5).Value = 588.2999877929688
$ws.Cells.Item(1333, 6).Value = "High"
$ws.Cells.Item(1333, 7).Value = 632.7999877929688
$ws.Cells.Item(1333, 8).Value = "week"
$ws.Cells.Item(1333, 9).Value = "17-06-2024 00:00:00"
$ws.Cells.Item(1333, 10).Value = 640
$ws.Cells.Item(1333, 11).Value = 629.7000122070312
$ws.Cells.Item(1333, 12).Value = "18/06/2024 11:36:33"

$ws.Cells.Item(1334, 1).Value = "INDHOTEL.NS"
$ws.Cells.Item(1334, 2).Value = 45383
$ws.Cells.Item(1334, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1334, 3).Value = 622.5
$ws.Cells.Item(1334, 4).Value = 594.25
$ws.Cells.Item(1334, 5).Value = 613.2999877929688
$ws.Cells.Item(1334, 6).Value = "High"
$ws.Cells.Item(1334, 7).Value = 622.5
$ws.Cells.Item(1334, 8).Value = "week"
$ws.Cells.Item(1334, 9).Value = "17-06-2024 00:00:00"
$ws.Cells.Item(1334, 10).Value = 626
$ws.Cells.Item(1334, 11).Value = 615.5999755859375
$ws.Cells.Item(1334, 12).Value = "18/06/2024 11:36:33"

$ws.Cells.Item(1335, 1).Value = "FACT.NS"
$ws.Cells.Item(1335, 2).Value = 45271
$ws.Cells.Item(1335, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1335, 3).Value = 861.1500244140625
$ws.Cells.Item(1335, 4).Value = 760.5
$ws.Cells.Item(1335, 5).Value = 801.0499877929688
$ws.Cells.Item(1335, 6).Value = "High"
$ws.Cells.Item(1335, 7).Value = 861.1500244140625
$ws.Cells.Item(1335, 8).Value = "week"
$ws.Cells.Item(1335, 9).Value = "17-06-2024 00:00:00"
$ws.Cells.Item(1335, 10).Value = 885
$ws.Cells.Item(1335, 11).Value = 810
$ws.Cells.Item(1335, 12).Value = "18/06/2024 11:36:33"

$ws.Cells.Item(1336, 1).Value = "POLICYBZR.NS"
$ws.Cells.Item(1336, 2).Value = 45390
$ws.Cells.Item(1336, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1336, 3).Value = 1400.349975585938
$ws.Cells.Item(1336, 4).Value = 1256.099975585938
$ws.Cells.Item(1336, 5).Value = 1264.150024414062
$ws.Cells.Item(1336, 6).Value = "High"
$ws.Cells.Item(1336, 7).Value = 1400.349975585938
$ws.Cells.Item(1336, 8).Value = "week"
$ws.Cells.Item(1336, 9).Value = "17-06-2024 00:00:00"
$ws.Cells.Item(1336, 10).Value = 1404.5
$ws.Cells.Item(1336, 11).Value = 1378.699951171875
$ws.Cells.Item(1336, 12).Value = "18/06/2024 11:36:33"

$ws.Cells.Item(1337, 1).Value = "SURANASOL.NS"
$ws.Cells.Item(1337, 2).Value = 45271
$ws.Cells.Item(1337, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1337, 3).Value = 39.90000152587891
$ws.Cells.Item(1337, 4).Value = 26.5
$ws.Cells.Item(1337, 5).Value = 35.54999923706055
$ws.Cells.Item(1337, 6).Value = "High"
$ws.Cells.Item(1337, 7).Value = 39.90000152587891
$ws.Cells.Item(1337, 8).Value = "week"
$ws.Cells.Item(1337, 9).Value = "17-06-2024 00:00:00"
$ws.Cells.Item(1337, 10).Value = 40
$ws.Cells.Item(1337, 11).Value = 39.45000076293945
$ws.Cells.Item(1337, 12).Value = "18/06/2024 11:36:33"

$ws.Cells.Item(1338, 1).Value = "CROWN.NS"
$ws.Cells.Item(1338, 2).Value = 45418
$ws.Cells.Item(1338, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1338, 3).Value = 261.6499938964844
$ws.Cells.Item(1338, 4).Value = 220.1999969482422
$ws.Cells.Item(1338, 5).Value = 232.0500030517578
$ws.Cells.Item(1338, 6).Value = "High"
$ws.Cells.Item(1338, 7).Value = 261.6499938964844
$ws.Cells.Item(1338, 8).Value = "week"
$ws.Cells.Item(1338, 9).Value = "17-06-2024 00:00:00"
$ws.Cells.Item(1338, 10).Value = 262.4800109863281
$ws.Cells.Item(1338, 11).Value = 249.9900054931641
$ws.Cells.Item(1338, 12).Value = "18/06/2024 11:36:33"

$ws.Cells.Item(1339, 1).Value = "BLUECLOUDS.BO"
$ws.Cells.Item(1339, 2).Value = 44158
$ws.Cells.Item(1339, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1339, 3).Value = 12
$ws.Cells.Item(1339, 4).Value = 12
$ws.Cells.Item(1339, 5).Value = 12
$ws.Cells.Item(1339, 6).Value = "High"
$ws.Cells.Item(1339, 7).Value = 12
$ws.Cells.Item(1339, 8).Value = "week"
$ws.Cells.Item(1339, 9).Value = "17-06-2024 00:00:00"
$ws.Cells.Item(1339, 10).Value = 151.5399932861328
$ws.Cells.Item(1339, 11).Value = 9.989999771118164
$ws.Cells.Item(1339, 12).Value = "18/06/2024 11:36:33"

$ws.Cells.Item(1340, 1).Value = "BLUECLOUDS.BO"
$ws.Cells.Item(1340, 2).Value = 44200
$ws.Cells.Item(1340, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1340, 3).Value = 12.52000045776367
$ws.Cells.Item(1340, 4).Value = 10.82999992370605
$ws.Cells.Item(1340, 5).Value = 11.35999965667725
$ws.Cells.Item(1340, 6).Value = "High"
$ws.Cells.Item(1340, 7).Value = 12.52000045776367
$ws.Cells.Item(1340, 8).Value = "week"
$ws.Cells.Item(1340, 9).Value = "17-06-2024 00:00:00"
$ws.Cells.Item(1340, 10).Value = 151.5399932861328
$ws.Cells.Item(1340, 11).Value = 9.989999771118164
$ws.Cells.Item(1340, 12).Value = "18/06/2024 11:36:33"

$ws.Cells.Item(1341, 1).Value = "BLUECLOUDS.BO"
$ws.Cells.Item(1341, 2).Value = 44277
$ws.Cells.Item(1341, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1341, 3).Value = 10.97000026702881
$ws.Cells.Item(1341, 4).Value = 10
$ws.Cells.Item(1341, 5).Value = 10
$ws.Cells.Item(1341, 6).Value = "High"
$ws.Cells.Item(1341, 7).Value = 10.97000026702881
$ws.Cells.Item(1341, 8).Value = "week"
$ws.Cells.Item(1341, 9).Value = "17-06-2024 00:00:00"
$ws.Cells.Item(1341, 10).Value = 151.5399932861328
$ws.Cells.Item(1341, 11).Value = 9.989999771118164
$ws.Cells.Item(1341, 12).Value = "18/06/2024 11:36:33"

$ws.Cells.Item(1342, 1).Value = "M&MFIN.NS"
$ws.Cells.Item(1342, 2).Value = 45327
$ws.Cells.Item(1342, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1342, 3).Value = 302.8999938964844
$ws.Cells.Item(1342, 4).Value = 282.1000061035156
$ws.Cells.Item(1342, 5).Value = 288.8500061035156
$ws.Cells.Item(1342, 6).Value = "High"
$ws.Cells.Item(1342, 7).Value = 302.8999938964844
$ws.Cells.Item(1342, 8).Value = "week"
$ws.Cells.Item(1342, 9).Value = "17-06-2024 00:00:00"
$ws.Cells.Item(1342, 10).Value = 313
$ws.Cells.Item(1342, 11).Value = 300
$ws.Cells.Item(1342, 12).Value = "18/06/2024 11:36:33"

$ws.Cells.Item(1343, 1).Value = "M&MFIN.NS"
$ws.Cells.Item(1343, 2).Value = 45390
$ws.Cells.Item(1343, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1343, 3).Value = 308.5
$ws.Cells.Item(1343, 4).Value = 295.75
$ws.Cells.Item(1343, 5).Value = 303.4500122070312
$ws.Cells.Item(1343, 6).Value = "High"
$ws.Cells.Item(1343, 7).Value = 308.5
$ws.Cells.Item(1343, 8).Value = "week"
$ws.Cells.Item(1343, 9).Value = "17-06-2024 00:00:00"
$ws.Cells.Item(1343, 10).Value = 313
$ws.Cells.Item(1343, 11).Value = 300
$ws.Cells.Item(1343, 12).Value = "18/06/2024 11:36:33"

